$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date field text from
#    14/06/2022 -> 22/06/2022 everywhere it appears: once on the
#    slide master and once on every slide layout's Date placeholder.
# ---------------------------------------------------------------------
$oldDate = "14/06/2022"
$newDate = "22/06/2022"

$master = $p.SlideMaster

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $shp = $master.Shapes.Item($si)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 9 ("4. Algoritmo" / glossario dei layer): finish the
#    "Conv2D" bullet ("... estrarre le feature e attraverso ..." ->
#    "... estrarre le features e attraverso ...", splitting the run so
#    "le features " is its own run) and fix the "Flattern" typo to
#    "Flatten".
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$contentShape = $slide9.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

$full = $tr.Text
$pos = $full.IndexOf("le feature ")
if ($pos -ge 0) {
    $sub = $tr.Characters($pos + 1, 11)
    $sub.Text = "le features "
}

$full2 = $tr.Text
$pos2 = $full2.IndexOf("Flattern")
if ($pos2 -ge 0) {
    $sub2 = $tr.Characters($pos2 + 1, 8)
    $sub2.Text = "Flatten"
}
